$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$fmtSrc = $ws.Range("B2")  # style-0 reference cell used to restore formatting after forcing text entry

$ws.Range('D2').Value = '71.921.98'
$ws.Range('E2').Value = '  +4.91%  '
$ws.Range('D3').Value = '4.044.37'
$ws.Range('E3').Value = '  +4.86%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '540.06'
$fmtSrc.Copy()
$ws.Range('D5').PasteSpecial(-4122)
$ws.Range('E5').Value = '  +3.56%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '153.64'
$fmtSrc.Copy()
$ws.Range('D6').PasteSpecial(-4122)
$ws.Range('E6').Value = '  +9.28%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.694'
$fmtSrc.Copy()
$ws.Range('D7').PasteSpecial(-4122)
$ws.Range('E7').Value = '  +14.35%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.764'
$fmtSrc.Copy()
$ws.Range('D9').PasteSpecial(-4122)
$ws.Range('E9').Value = '  +7.46%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.175'
$fmtSrc.Copy()
$ws.Range('D10').PasteSpecial(-4122)
$ws.Range('E10').Value = '  +4.73%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000332'
$fmtSrc.Copy()
$ws.Range('D11').PasteSpecial(-4122)
$ws.Range('E11').Value = '  +3.57%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '48.21'
$fmtSrc.Copy()
$ws.Range('D12').PasteSpecial(-4122)
$ws.Range('E12').Value = '  +16.09%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.88'
$fmtSrc.Copy()
$ws.Range('D13').PasteSpecial(-4122)
$ws.Range('E13').Value = '  +4.30%  '
$ws.Range('D14').Value = '4.688.91'
$ws.Range('E14').Value = '  +4.82%  '
$ws.Range('D15').Value = '4.036.76'
$ws.Range('E15').Value = '  +3.44%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.40'
$fmtSrc.Copy()
$ws.Range('D16').PasteSpecial(-4122)
$ws.Range('E16').Value = '  +2.17%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '20.71'
$fmtSrc.Copy()
$ws.Range('D17').PasteSpecial(-4122)
$ws.Range('E17').Value = '  -2.84%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.22'
$fmtSrc.Copy()
$ws.Range('D18').PasteSpecial(-4122)
$ws.Range('E18').Value = '  +2.17%  '
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = '71.864.34'
$ws.Range('E20').Value = '  +4.79%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '435.62'
$fmtSrc.Copy()
$ws.Range('D21').PasteSpecial(-4122)
$ws.Range('E21').Value = '  +4.80%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '99.37'
$fmtSrc.Copy()
$ws.Range('D22').PasteSpecial(-4122)
$ws.Range('E22').Value = '  +14.44%  '
$ws.Range('E23').Value = '  +3.08%  '
$ws.Range('E24').Value = '  +7.08%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '14.75'
$fmtSrc.Copy()
$ws.Range('D25').PasteSpecial(-4122)
$ws.Range('E25').Value = '  +5.54%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.30'
$fmtSrc.Copy()
$ws.Range('D26').PasteSpecial(-4122)
$ws.Range('E26').Value = '  -2.94%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.02'
$fmtSrc.Copy()
$ws.Range('D27').PasteSpecial(-4122)
$ws.Range('E27').Value = '  +5.26%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '37.16'
$fmtSrc.Copy()
$ws.Range('D28').PasteSpecial(-4122)
$ws.Range('E28').Value = '  +4.71%  '
$ws.Range('E29').Value = '  +2.90%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.63'
$fmtSrc.Copy()
$ws.Range('D30').PasteSpecial(-4122)
$ws.Range('E30').Value = '  +30.41%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '13.67'
$fmtSrc.Copy()
$ws.Range('D31').PasteSpecial(-4122)
$ws.Range('E31').Value = '  +1.99%  '
$ws.Range('E32').Value = '  +5.75%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '691.02'
$fmtSrc.Copy()
$ws.Range('D33').PasteSpecial(-4122)
$ws.Range('E33').Value = '  +1.82%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.91'
$fmtSrc.Copy()
$ws.Range('D34').PasteSpecial(-4122)
$ws.Range('E34').Value = '  +1.26%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '67.31'
$fmtSrc.Copy()
$ws.Range('D35').PasteSpecial(-4122)
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '43.42'
$fmtSrc.Copy()
$ws.Range('D36').PasteSpecial(-4122)
$ws.Range('E36').Value = '  +10.12%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.437'
$fmtSrc.Copy()
$ws.Range('D37').PasteSpecial(-4122)
$ws.Range('E37').Value = '  -2.01%  '
$ws.Range('E38').Value = '  +6.36%  '
$ws.Range('D39').Value = '0.0₃0844'
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.46'
$fmtSrc.Copy()
$ws.Range('D40').PasteSpecial(-4122)
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.39'
$fmtSrc.Copy()
$ws.Range('D41').PasteSpecial(-4122)
$ws.Range('E41').Value = '  +6.67%  '
$ws.Range('B42').Value = 'Dai'
$ws.Range('C42').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.00'
$fmtSrc.Copy()
$ws.Range('D42').PasteSpecial(-4122)
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0496'
$fmtSrc.Copy()
$ws.Range('D43').PasteSpecial(-4122)
$ws.Range('E43').Value = '  +4.76%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.999'
$fmtSrc.Copy()
$ws.Range('D44').PasteSpecial(-4122)
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('E45').Value = '  +7.76%  '
$ws.Range('E46').Value = '  -3.80%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.42'
$fmtSrc.Copy()
$ws.Range('D47').PasteSpecial(-4122)
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.55'
$fmtSrc.Copy()
$ws.Range('D48').PasteSpecial(-4122)
$ws.Range('E48').Value = '  +10.21%  '
$ws.Range('E49').Value = '  +2.25%  '
$ws.Range('E50').Value = '  +2.65%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.000271'
$fmtSrc.Copy()
$ws.Range('D51').PasteSpecial(-4122)
$ws.Range('E51').Value = '  -0.78%  '

$excel.CutCopyMode = $false
